$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 9).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 13).ClearContents()
$ws.Cells.Item(9, 14).ClearContents()
$ws.Cells.Item(33, 8).Value = 139.25
$ws.Cells.Item(33, 9).Value = 144.6
$ws.Cells.Item(33, 11).Value = 144.6
$ws.Cells.Item(33, 13).Value = 84.40000000000001
$ws.Cells.Item(43, 8).Value = 13448.75
$ws.Cells.Item(43, 9).Value = 16931.666
$ws.Cells.Item(43, 10).Value = 3000
$ws.Cells.Item(43, 11).Value = 16931.666
$ws.Cells.Item(43, 12).Value = 3000
$ws.Cells.Item(43, 13).Value = -16862.666
$ws.Cells.Item(43, 14).Value = -3138
$ws.Cells.Item(74, 8).Value = 4957.4194
$ws.Cells.Item(74, 9).Value = 4760.8
$ws.Cells.Item(74, 10).Value = 5776.6665
$ws.Cells.Item(74, 11).Value = 4760.8
$ws.Cells.Item(74, 12).Value = 5776.6665
$ws.Cells.Item(74, 13).Value = -3824.8
$ws.Cells.Item(74, 14).Value = -7648.6665
$ws.Cells.Item(76, 8).Value = 5720.5713
$ws.Cells.Item(76, 9).Value = 3655.4285
$ws.Cells.Item(76, 10).Value = 7785.7144
$ws.Cells.Item(76, 11).Value = 3655.4285
$ws.Cells.Item(76, 12).Value = 7785.7144
$ws.Cells.Item(76, 13).Value = -3340.4285
$ws.Cells.Item(76, 14).Value = -8415.714400000001
$ws.Cells.Item(77, 8).Value = 4957.4194
$ws.Cells.Item(77, 9).Value = 4760.8
$ws.Cells.Item(77, 10).Value = 5776.6665
$ws.Cells.Item(77, 11).Value = 23804
$ws.Cells.Item(77, 12).Value = 28883.3325
$ws.Cells.Item(77, 13).Value = -19124
$ws.Cells.Item(77, 14).Value = -38243.3325
$ws.Cells.Item(79, 8).Value = 5720.5713
$ws.Cells.Item(79, 9).Value = 3655.4285
$ws.Cells.Item(79, 10).Value = 7785.7144
$ws.Cells.Item(79, 11).Value = 3655.4285
$ws.Cells.Item(79, 12).Value = 7785.7144
$ws.Cells.Item(79, 13).Value = -2563.4285
$ws.Cells.Item(79, 14).Value = -9969.714400000001
$ws.Cells.Item(80, 8).Value = 555.36365
$ws.Cells.Item(80, 9).Value = 473
$ws.Cells.Item(80, 10).Value = 699.5
$ws.Cells.Item(80, 11).Value = 1419
$ws.Cells.Item(80, 12).Value = 2098.5
$ws.Cells.Item(80, 13).Value = -421
$ws.Cells.Item(80, 14).Value = -4094.5
$ws.Cells.Item(83, 8).Value = 555.36365
$ws.Cells.Item(83, 9).Value = 473
$ws.Cells.Item(83, 10).Value = 699.5
$ws.Cells.Item(83, 11).Value = 4257
$ws.Cells.Item(83, 12).Value = 6295.5
$ws.Cells.Item(83, 13).Value = 735
$ws.Cells.Item(83, 14).Value = -16279.5
$ws.Cells.Item(92, 8).Value = 616.5333000000001
$ws.Cells.Item(92, 9).Value = 618.1818
$ws.Cells.Item(92, 11).Value = 618.1818
$ws.Cells.Item(92, 13).Value = 629.8182
$ws.Cells.Item(99, 8).Value = 0
$ws.Cells.Item(99, 9).Value = 0
$ws.Cells.Item(99, 11).Value = 0
$ws.Cells.Item(99, 13).ClearContents()
$ws.Cells.Item(118, 8).Value = 224.55556
$ws.Cells.Item(118, 9).Value = 194.4
$ws.Cells.Item(118, 10).Value = 262.25
$ws.Cells.Item(118, 11).Value = 583.2
$ws.Cells.Item(118, 12).Value = 786.75
$ws.Cells.Item(118, 13).Value = 1073.8
$ws.Cells.Item(118, 14).Value = -4100.75
$ws.Cells.Item(137, 8).Value = 1624.75
$ws.Cells.Item(137, 9).Value = 1624.75
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 4874.25
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -2324.25
$ws.Cells.Item(137, 14).ClearContents()
$ws.Cells.Item(138, 8).Value = 2132.4194
$ws.Cells.Item(138, 10).Value = 2994.5334
$ws.Cells.Item(138, 12).Value = 8983.600199999999
$ws.Cells.Item(138, 14).Value = -19263.6002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 649.8823
$ws.Cells.Item(97, 9).Value = 565.5
$ws.Cells.Item(97, 10).Value = 2000
$ws.Cells.Item(97, 11).Value = 565.5
$ws.Cells.Item(97, 12).Value = 2000
$ws.Cells.Item(97, 13).Value = -69.5
$ws.Cells.Item(97, 14).Value = -2992

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 26660
$ws.Cells.Item(81, 10).Value = 19990
$ws.Cells.Item(81, 12).Value = 19990
$ws.Cells.Item(81, 14).Value = -22112
$ws.Cells.Item(84, 8).Value = 26660
$ws.Cells.Item(84, 10).Value = 19990
$ws.Cells.Item(84, 12).Value = 59970
$ws.Cells.Item(84, 14).Value = -70578
$ws.Cells.Item(94, 8).Value = 2350
$ws.Cells.Item(94, 9).Value = 1820
$ws.Cells.Item(94, 10).Value = 5000
$ws.Cells.Item(94, 11).Value = 1820
$ws.Cells.Item(94, 12).Value = 5000
$ws.Cells.Item(94, 13).Value = -1369
$ws.Cells.Item(94, 14).Value = -5902

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4096.7856
$ws.Cells.Item(31, 9).Value = 2734.5
$ws.Cells.Item(31, 11).Value = 2734.5
$ws.Cells.Item(31, 13).Value = -2439.5
$ws.Cells.Item(34, 8).Value = 4096.7856
$ws.Cells.Item(34, 9).Value = 2734.5
$ws.Cells.Item(34, 11).Value = 2734.5
$ws.Cells.Item(34, 13).Value = -2532.5
$ws.Cells.Item(122, 8).Value = 4147.1665
$ws.Cells.Item(122, 9).Value = 3296.3333
$ws.Cells.Item(122, 10).Value = 4998
$ws.Cells.Item(122, 11).Value = 9888.999899999999
$ws.Cells.Item(122, 12).Value = 14994
$ws.Cells.Item(122, 13).Value = -7438.999899999999
$ws.Cells.Item(122, 14).Value = -19894
$ws.Cells.Item(134, 8).Value = 2719.5
$ws.Cells.Item(134, 9).Value = 2399.75
$ws.Cells.Item(134, 11).Value = 7199.25
$ws.Cells.Item(134, 13).Value = -4664.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value = 1940.2368
$ws.Cells.Item(11, 9).Value = 1792.5555
$ws.Cells.Item(11, 10).Value = 1986.069
$ws.Cells.Item(11, 11).Value = 5377.666499999999
$ws.Cells.Item(11, 12).Value = 5958.207
$ws.Cells.Item(11, 13).Value = -5237.666499999999
$ws.Cells.Item(11, 14).Value = -6238.207
$ws.Cells.Item(131, 8).Value = 2444.0754
$ws.Cells.Item(131, 9).Value = 1381
$ws.Cells.Item(131, 10).Value = 2507.86
$ws.Cells.Item(131, 11).Value = 4143
$ws.Cells.Item(131, 12).Value = 7523.58
$ws.Cells.Item(131, 13).Value = 897
$ws.Cells.Item(131, 14).Value = -17603.58

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3098.9092
$ws.Cells.Item(102, 9).Value = 1147
$ws.Cells.Item(102, 10).Value = 4214.2856
$ws.Cells.Item(102, 11).Value = 1147
$ws.Cells.Item(102, 12).Value = 4214.2856
$ws.Cells.Item(102, 13).Value = 475
$ws.Cells.Item(102, 14).Value = -7458.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(29, 8).Value = 29999.5
$ws.Cells.Item(29, 9).Value = 29999.5
$ws.Cells.Item(29, 11).Value = 29999.5
$ws.Cells.Item(29, 13).Value = -29704.5
$ws.Cells.Item(35, 8).Value = 6500
$ws.Cells.Item(35, 9).Value = 6500
$ws.Cells.Item(35, 10).Value = 0
$ws.Cells.Item(35, 11).Value = 6500
$ws.Cells.Item(35, 12).Value = 0
$ws.Cells.Item(35, 13).Value = -6164
$ws.Cells.Item(35, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 1000004
$ws.Cells.Item(132, 9).Value = 1000004
$ws.Cells.Item(132, 11).Value = 3000012
$ws.Cells.Item(132, 13).Value = -2997482

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 2749.5
$ws.Cells.Item(32, 9).Value = 1500
$ws.Cells.Item(32, 10).Value = 3999
$ws.Cells.Item(32, 11).Value = 1500
$ws.Cells.Item(32, 12).Value = 3999
$ws.Cells.Item(32, 13).Value = -1183
$ws.Cells.Item(32, 14).Value = -4633
